# Update the crypto price/volume table (cols D and E) with refreshed
# values, as produced by the scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = '@'
$dCell.Value = '58.144.71'
$eCell = $ws.Range("E2")
$eCell.NumberFormat = '@'
$eCell.Value = '  -1.30%  '

$dCell = $ws.Range("D3")
$dCell.NumberFormat = '@'
$dCell.Value = '3.116.72'
$eCell = $ws.Range("E3")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.76%  '

$dCell = $ws.Range("D4")
$dCell.NumberFormat = '@'
$dCell.Value = '1.00'
$eCell = $ws.Range("E4")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.06%  '

$dCell = $ws.Range("D5")
$dCell.NumberFormat = '@'
$dCell.Value = '527.66'
$eCell = $ws.Range("E5")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.98%  '

$dCell = $ws.Range("D6")
$dCell.NumberFormat = '@'
$dCell.Value = '142.43'
$eCell = $ws.Range("E6")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.83%  '

$dCell = $ws.Range("D7")
$dCell.NumberFormat = '@'
$dCell.Value = '1.00'
$eCell = $ws.Range("E7")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.09%  '

$dCell = $ws.Range("D8")
$dCell.NumberFormat = '@'
$dCell.Value = '3.115.60'
$eCell = $ws.Range("E8")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.76%  '

$dCell = $ws.Range("D9")
$dCell.NumberFormat = '@'
$dCell.Value = '0.445'
$eCell = $ws.Range("E9")
$eCell.NumberFormat = '@'
$eCell.Value = '  +1.41%  '

$dCell = $ws.Range("D10")
$dCell.NumberFormat = '@'
$dCell.Value = '7.17'
$eCell = $ws.Range("E10")
$eCell.NumberFormat = '@'
$eCell.Value = '  -2.66%  '

$dCell = $ws.Range("D11")
$dCell.NumberFormat = '@'
$dCell.Value = '0.109'
$eCell = $ws.Range("E11")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.42%  '

$dCell = $ws.Range("D12")
$dCell.NumberFormat = '@'
$dCell.Value = '0.393'
$eCell = $ws.Range("E12")
$eCell.NumberFormat = '@'
$eCell.Value = '  +2.89%  '

$dCell = $ws.Range("D13")
$dCell.NumberFormat = '@'
$dCell.Value = '3.651.53'
$eCell = $ws.Range("E13")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.59%  '

$dCell = $ws.Range("D14")
$dCell.NumberFormat = '@'
$dCell.Value = '0.134'
$eCell = $ws.Range("E14")
$eCell.NumberFormat = '@'
$eCell.Value = '  +3.44%  '

$dCell = $ws.Range("D15")
$dCell.NumberFormat = '@'
$dCell.Value = '25.84'
$eCell = $ws.Range("E15")
$eCell.NumberFormat = '@'
$eCell.Value = '  -3.77%  '

$dCell = $ws.Range("D16")
$dCell.NumberFormat = '@'
$dCell.Value = '0.0000166'
$eCell = $ws.Range("E16")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.31%  '

$dCell = $ws.Range("D17")
$dCell.NumberFormat = '@'
$dCell.Value = '58.177.75'
$eCell = $ws.Range("E17")
$eCell.NumberFormat = '@'
$eCell.Value = '  -1.34%  '

$dCell = $ws.Range("D18")
$dCell.NumberFormat = '@'
$dCell.Value = '3.112.32'
$eCell = $ws.Range("E18")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.52%  '

$dCell = $ws.Range("D19")
$dCell.NumberFormat = '@'
$dCell.Value = '6.14'
$eCell = $ws.Range("E19")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.24%  '

$dCell = $ws.Range("D20")
$dCell.NumberFormat = '@'
$dCell.Value = '12.77'
$eCell = $ws.Range("E20")
$eCell.NumberFormat = '@'
$eCell.Value = '  -1.24%  '

$dCell = $ws.Range("D21")
$dCell.NumberFormat = '@'
$dCell.Value = '8.00'
$eCell = $ws.Range("E21")
$eCell.NumberFormat = '@'
$eCell.Value = '  -1.43%  '

$dCell = $ws.Range("D22")
$dCell.NumberFormat = '@'
$dCell.Value = '343.42'
$eCell = $ws.Range("E22")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.03%  '

$dCell = $ws.Range("D23")
$dCell.NumberFormat = '@'
$dCell.Value = '1.00'
$eCell = $ws.Range("E23")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.04%  '

$dCell = $ws.Range("D24")
$dCell.NumberFormat = '@'
$dCell.Value = '0.516'
$eCell = $ws.Range("E24")
$eCell.NumberFormat = '@'
$eCell.Value = '  +1.99%  '

$dCell = $ws.Range("D25")
$dCell.NumberFormat = '@'
$dCell.Value = '67.64'
$eCell = $ws.Range("E25")
$eCell.NumberFormat = '@'
$eCell.Value = '  +2.94%  '

$dCell = $ws.Range("D26")
$dCell.NumberFormat = '@'
$dCell.Value = '0.170'
$eCell = $ws.Range("E26")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.78%  '

$dCell = $ws.Range("D27")
$dCell.NumberFormat = '@'
$dCell.Value = '0.999'
$eCell = $ws.Range("E27")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.15%  '

$dCell = $ws.Range("D28")
$dCell.NumberFormat = '@'
$dCell.Value = '0.0₃0929'
$eCell = $ws.Range("E28")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.73%  '

$dCell = $ws.Range("D29")
$dCell.NumberFormat = '@'
$dCell.Value = '1.00'
$eCell = $ws.Range("E29")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.11%  '

$dCell = $ws.Range("D30")
$dCell.NumberFormat = '@'
$dCell.Value = '6.39'
$eCell = $ws.Range("E30")
$eCell.NumberFormat = '@'
$eCell.Value = '  -4.55%  '

$dCell = $ws.Range("D31")
$dCell.NumberFormat = '@'
$dCell.Value = '7.29'
$eCell = $ws.Range("E31")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.65%  '

$dCell = $ws.Range("D32")
$dCell.NumberFormat = '@'
$dCell.Value = '1.88'
$eCell = $ws.Range("E32")
$eCell.NumberFormat = '@'
$eCell.Value = '  +2.00%  '

$dCell = $ws.Range("D33")
$dCell.NumberFormat = '@'
$dCell.Value = '21.09'
$eCell = $ws.Range("E33")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.61%  '

$dCell = $ws.Range("D34")
$dCell.NumberFormat = '@'
$dCell.Value = '1.20'
$eCell = $ws.Range("E34")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.90%  '

$dCell = $ws.Range("D35")
$dCell.NumberFormat = '@'
$dCell.Value = '158.46'
$eCell = $ws.Range("E35")
$eCell.NumberFormat = '@'
$eCell.Value = '  +2.32%  '

$dCell = $ws.Range("D36")
$dCell.NumberFormat = '@'
$dCell.Value = '4.65'
$eCell = $ws.Range("E36")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.48%  '

$dCell = $ws.Range("D37")
$dCell.NumberFormat = '@'
$dCell.Value = '6.20'
$eCell = $ws.Range("E37")
$eCell.NumberFormat = '@'
$eCell.Value = '  +1.03%  '

$dCell = $ws.Range("D38")
$dCell.NumberFormat = '@'
$dCell.Value = '26.43'
$eCell = $ws.Range("E38")
$eCell.NumberFormat = '@'
$eCell.Value = '  -1.83%  '

$dCell = $ws.Range("D39")
$dCell.NumberFormat = '@'
$dCell.Value = '1.25'
$eCell = $ws.Range("E39")
$eCell.NumberFormat = '@'
$eCell.Value = '  -3.10%  '

$dCell = $ws.Range("D40")
$dCell.NumberFormat = '@'
$dCell.Value = '1.65'
$eCell = $ws.Range("E40")
$eCell.NumberFormat = '@'
$eCell.Value = '  +13.33%  '

$dCell = $ws.Range("D41")
$dCell.NumberFormat = '@'
$dCell.Value = '0.0670'
$eCell = $ws.Range("E41")
$eCell.NumberFormat = '@'
$eCell.Value = '  -2.24%  '

$dCell = $ws.Range("D42")
$dCell.NumberFormat = '@'
$dCell.Value = '4.00'
$eCell = $ws.Range("E42")
$eCell.NumberFormat = '@'
$eCell.Value = '  +1.72%  '

$dCell = $ws.Range("D43")
$dCell.NumberFormat = '@'
$dCell.Value = '0.690'
$eCell = $ws.Range("E43")
$eCell.NumberFormat = '@'
$eCell.Value = '  +3.72%  '

$dCell = $ws.Range("D44")
$dCell.NumberFormat = '@'
$dCell.Value = '3.157.28'
$eCell = $ws.Range("E44")
$eCell.NumberFormat = '@'
$eCell.Value = '  +0.63%  '

$dCell = $ws.Range("D45")
$dCell.NumberFormat = '@'
$dCell.Value = '36.52'
$eCell = $ws.Range("E45")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.75%  '

$dCell = $ws.Range("D46")
$dCell.NumberFormat = '@'
$dCell.Value = '0.999'
$eCell = $ws.Range("E46")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.12%  '

$dCell = $ws.Range("D47")
$dCell.NumberFormat = '@'
$dCell.Value = '0.0263'
$eCell = $ws.Range("E47")
$eCell.NumberFormat = '@'
$eCell.Value = '  +3.00%  '

$dCell = $ws.Range("D48")
$dCell.NumberFormat = '@'
$dCell.Value = '2.273.04'
$eCell = $ws.Range("E48")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.49%  '

$dCell = $ws.Range("D49")
$dCell.NumberFormat = '@'
$dCell.Value = '0.999'
$eCell = $ws.Range("E49")
$eCell.NumberFormat = '@'
$eCell.Value = '  +3.94%  '

$dCell = $ws.Range("D50")
$dCell.NumberFormat = '@'
$dCell.Value = '6.13'
$eCell = $ws.Range("E50")
$eCell.NumberFormat = '@'
$eCell.Value = '  +2.23%  '

$dCell = $ws.Range("D51")
$dCell.NumberFormat = '@'
$dCell.Value = '20.66'
$eCell = $ws.Range("E51")
$eCell.NumberFormat = '@'
$eCell.Value = '  -0.95%  '

